$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (matches source data format)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.744.21'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.628.26'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '215.07'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.2559'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '0.06317'
$ws.Range("E9").Value = '  -0.44%  '
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("D11").Value = '0.07771'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '1.644.98'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = '4.218'
$ws.Range("E13").Value = '  -1.59%  '
$ws.Range("D14").Value = '1.848.69'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").Value = '0.5497'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").Value = '0.0₅7535'
$ws.Range("E17").Value = '  -2.18%  '
$ws.Range("D18").Value = '25.779.64'
$ws.Range("E18").Value = '  -1.00%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '4.393'
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("D21").Value = '193.19'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("D22").Value = '9.820'
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("D23").Value = '5.979'
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").Value = '1.888'
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("D26").Value = '141.99'
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("E27").Value = '  +5.26%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '15.52'
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '6.717'
$ws.Range("E29").Value = '  -1.30%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").Value = '0.04864'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").Value = '3.222'
$ws.Range("E32").Value = '  -0.97%  '
$ws.Range("D33").Value = '3.144'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").Value = '1.534'
$ws.Range("E34").Value = '  +0.65%  '
$ws.Range("D35").Value = '2.382'
$ws.Range("E35").Value = '  +0.70%  '
$ws.Range("D36").Value = '0.8905'
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("D37").Value = '2.531'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("D38").Value = '0.5484'
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = '1.108.34'
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").Value = '0.01544'
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '5.532'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("D43").Value = '0.7959'
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("D44").Value = '97.06'
$ws.Range("E44").Value = '  -2.30%  '
$ws.Range("D45").Value = '1.772.82'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E46").Value = '  -12.97%  '
$ws.Range("D47").Value = '0.4429'
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").Value = '54.45'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").Value = '0.05133'
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").Value = '7.501'
$ws.Range("E51").Value = '  +2.66%  '
